$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.386.67"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "2.216.93"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "109.14"
$ws.Range("E5").Value = "  -11.89%  "
$ws.Range("D6").Value = "297.97"
$ws.Range("E6").Value = "  +11.79%  "
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").Value = "44.78"
$ws.Range("E10").Value = "  -7.68%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").Value = "54.40"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "8.80"
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "0.949"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "15.04"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("D17").Value = "2.551.29"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "2.239.15"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "42.325.10"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("D22").Value = "73.90"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").Value = "3.45"
$ws.Range("E23").Value = "  +19.56%  "
$ws.Range("E24").Value = "  -6.71%  "
$ws.Range("D25").Value = "229.27"
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("D26").Value = "9.27"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").Value = "11.65"
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("D29").Value = "3.89"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "38.41"
$ws.Range("E31").Value = "  -10.48%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("D33").Value = "174.44"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").Value = "21.04"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").Value = "0.0881"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "5.66"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").Value = "4.84"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("D38").Value = "4.24"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.126"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").Value = "0.0366"
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").Value = "70.73"
$ws.Range("E44").Value = "  -6.12%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "12.71"
$ws.Range("E46").Value = "  -9.31%  "
$ws.Range("E47").Value = "  -4.74%  "
$ws.Range("D48").Value = "5.46"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("D50").Value = "103.22"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").Value = "8.50"
$ws.Range("E51").Value = "  -1.35%  "
